$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4,5,6,7,8,9,10,11,12,14,21,28,37,38,39,63)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

$ws.Activate()
$ws.Range("E31").Select() | Out-Null
$excel.ActiveWindow.Zoom = 70
